# Extended forward-stationing logic for multiple forward demands
# and added new test data to SupplyDemand and the timeline.
#
# Adds a new "UtahFwd" row (row 7) to the single timeline worksheet with
# the same StartDay/Duration values as the existing UtahFwd-style rows
# (ForceCode=UtahFwd, StartDay=91, Duration=273), mirroring the AlaskaFwd
# row already present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "UtahFwd"
$ws.Range("B7").Value = 91
$ws.Range("C7").Value = 273

# Move the active selection, matching the author's final cursor position.
$ws.Range("G21").Select()
